$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates ---
$ws.Range("G3").Value = "Visualizza lo stato delle aste in corso di un utente"
$ws.Range("G7").Value = "Imposta Controfferta"
$ws.Range("D12").Value = "Riceve"

# --- Sheet view change: scroll so column F is the left-most visible column ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

# --- Formula updates ---
$ws.Range("J3").Formula = "=H3*(1+1+1+1+(B2*B3/2)*2)"
$ws.Range("E4").Formula = "=B2*B3*(E8+E9)"
$ws.Range("H4").Formula = "=H10/B11*B3"
$ws.Range("J4").Formula = "=H4*((1+1+1)*2+1)"
$ws.Range("H5").Formula = "=H10"
$ws.Range("J5").Formula = "=H5*((B5+E12*2+E9*4))"
$ws.Range("H6").Formula = "=100"
$ws.Range("J6").Formula = "=H6*((E14-E12)*2+E8*2)"
$ws.Range("H7").Formula = "=H4*0.1"
$ws.Range("J7").Formula = "=H7*2*2"
$ws.Range("E11").Formula = "=E4"
$ws.Range("E14").Formula = "=E4"

# --- Page setup: clear custom first page number ---
$ws.PageSetup.FirstPageNumber = 1
